$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.899.21"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.16%  "

$ws.Range("D3").Value = "'1.879.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "'324.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").Value = "'0.4619"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  -2.33%  "

$ws.Range("D9").Value = "'0.07859"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").Value = "'0.9841"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.09%  "

$ws.Range("E11").Value = "  -2.35%  "

$ws.Range("D12").Value = "'1.819.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.93%  "

$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").Value = "'5.667"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("D15").Value = "'0.06974"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("D16").Value = "'88.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.00%  "

$ws.Range("D18").Value = "'0.000009943"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.82%  "

$ws.Range("D19").Value = "'16.97"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.41%  "

$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").Value = "'28.887.80"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.23%  "

$ws.Range("D22").Value = "'5.267"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.63%  "

$ws.Range("E23").Value = "  -1.77%  "

$ws.Range("D24").Value = "'2.104"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("D25").Value = "'156.25"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").Value = "'19.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("D27").Value = "'5.909"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("D28").Value = "'117.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.80%  "

$ws.Range("E29").Value = "  -7.33%  "

$ws.Range("D30").Value = "'0.09367"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").Value = "'0.9022"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.97%  "

$ws.Range("D32").Value = "'5.275"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("D34").Value = "'3.249"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("B35").Value = "TrustWalletToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D35").Value = "'1.172"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.66%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.05746"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.46%  "

$ws.Range("D37").Value = "'0.02076"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("D38").Value = "'1.000"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5661"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.37%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.622"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.83%  "

$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("D42").Value = "'9.689"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.16%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'11.97"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.238"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.47%  "

$ws.Range("D45").Value = "'0.5338"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.58%  "

$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("E47").Value = "  -4.68%  "

$ws.Range("D48").Value = "'2.538"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").Value = "'112.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.03%  "

$ws.Range("D50").Value = "'1.066"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.78%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'70.68"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.85%  "
